# RDCC-5182 Added Version check
# Adds a new "VERSION" worksheet at the end of the workbook containing a
# "File version" / "vx.xx" label pair, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet (so it lands at the end,
# as the tab order, rather than Excel's default of inserting before the
# active sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws.Name = "VERSION"

# Populate the version info.
$ws.Range("A6").Value = "File version"
$ws.Range("B6").Value = "vx.xx"

# Make the new sheet the active tab/selection, matching the committed file.
$ws.Range("B6").Select()
$ws.Activate()
